# Auto-generated Excel COM-interop script to update cryptos price table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.747.03"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.724.63"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.46"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.02"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.108"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.374"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.205.81"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.588.62"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.726.37"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.23"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.69"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.04"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.50"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.512"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.32"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.26"
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0897"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +8.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.14"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.50"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.88"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.93"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "349.83"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.961"
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.22"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.44"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.51"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.76"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0577"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0248"
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.22"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.12%  "
